$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 2).Value2 = 12.68435830011695
$ws.Cells.Item(2, 3).Value2 = 7.0629453439059
$ws.Cells.Item(2, 4).Value2 = 6.010429346583905
$ws.Cells.Item(2, 5).Value2 = 11.72728366628107
$ws.Cells.Item(2, 7).Value2 = 3.676908119126252
$ws.Cells.Item(2, 9).Value2 = 26.87708286153848
$ws.Cells.Item(2, 11).Value2 = 10.25798324992586
$ws.Cells.Item(2, 12).Value2 = 10.21497529931431
$ws.Cells.Item(2, 13).Value2 = 14.42688439708234
$ws.Cells.Item(2, 14).Value2 = 20.88460521208003
$ws.Cells.Item(2, 15).Value2 = 27.53448474598377

$ws.Cells.Item(3, 2).Value2 = 12.48407508107757
$ws.Cells.Item(3, 3).Value2 = 6.980345290816367
$ws.Cells.Item(3, 4).Value2 = 5.895938062469471
$ws.Cells.Item(3, 5).Value2 = 11.7505960917783
$ws.Cells.Item(3, 7).Value2 = 3.67882824378863
$ws.Cells.Item(3, 9).Value2 = 26.95276870590538
$ws.Cells.Item(3, 11).Value2 = 10.11136288809548
$ws.Cells.Item(3, 12).Value2 = 10.22305001863564
$ws.Cells.Item(3, 13).Value2 = 14.40153587745712
$ws.Cells.Item(3, 14).Value2 = 20.94537179274289
$ws.Cells.Item(3, 15).Value2 = 27.59573036480757

$ws.Cells.Item(4, 2).Value2 = 12.36205118020513
$ws.Cells.Item(4, 3).Value2 = 6.928283206292749
$ws.Cells.Item(4, 4).Value2 = 5.826257488345884
$ws.Cells.Item(4, 5).Value2 = 11.76625441815842
$ws.Cells.Item(4, 7).Value2 = 3.680070629292825
$ws.Cells.Item(4, 9).Value2 = 27.00363282934341
$ws.Cells.Item(4, 11).Value2 = 10.0222259354613
$ws.Cells.Item(4, 12).Value2 = 10.22940318338174
$ws.Cells.Item(4, 13).Value2 = 14.38814989319113
$ws.Cells.Item(4, 14).Value2 = 20.98444049382451
$ws.Cells.Item(4, 15).Value2 = 27.638425952848

$ws.Cells.Item(5, 2).Value2 = 12.31262916874295
$ws.Cells.Item(5, 3).Value2 = 6.906738527786207
$ws.Cells.Item(5, 4).Value2 = 5.798060766411136
$ws.Cells.Item(5, 5).Value2 = 11.77297375526686
$ws.Cells.Item(5, 7).Value2 = 3.680592908873351
$ws.Cells.Item(5, 9).Value2 = 27.02546426856997
$ws.Cells.Item(5, 11).Value2 = 9.986169672225667
$ws.Cells.Item(5, 12).Value2 = 10.23234352259543
$ws.Cells.Item(5, 13).Value2 = 14.38324706298273
$ws.Cells.Item(5, 14).Value2 = 21.00080456520439
$ws.Cells.Item(5, 15).Value2 = 27.65710326363788

$ws.Cells.Item(6, 2).Value2 = 12.30444291453196
$ws.Cells.Item(6, 3).Value2 = 6.903141459189125
$ws.Cells.Item(6, 4).Value2 = 5.793391994889914
$ws.Cells.Item(6, 5).Value2 = 11.77410994815562
$ws.Cells.Item(6, 7).Value2 = 3.680680600663902
$ws.Cells.Item(6, 9).Value2 = 27.0291560143087
$ws.Cells.Item(6, 11).Value2 = 9.980200020765952
$ws.Cells.Item(6, 12).Value2 = 10.23285299974707
$ws.Cells.Item(6, 13).Value2 = 14.38246641564698
$ws.Cells.Item(6, 14).Value2 = 21.00354861535098
$ws.Cells.Item(6, 15).Value2 = 27.66028178755862

$ws.Cells.Item(7, 2).Value2 = 12.36138334251216
$ws.Cells.Item(7, 3).Value2 = 6.927993966660911
$ws.Cells.Item(7, 4).Value2 = 5.82587635551483
$ws.Cells.Item(7, 5).Value2 = 11.76634366664428
$ws.Cells.Item(7, 7).Value2 = 3.68007760809875
$ws.Cells.Item(7, 9).Value2 = 27.00392278707978
$ws.Cells.Item(7, 11).Value2 = 10.02173852549828
$ws.Cells.Item(7, 12).Value2 = 10.2294414144835
$ws.Cells.Item(7, 13).Value2 = 14.38808153111574
$ws.Cells.Item(7, 14).Value2 = 20.98465938902253
$ws.Cells.Item(7, 15).Value2 = 27.63867266712376

$ws.Cells.Item(8, 2).Value2 = 12.61514203282658
$ws.Cells.Item(8, 3).Value2 = 7.034749591185789
$ws.Cells.Item(8, 4).Value2 = 5.970850504865891
$ws.Cells.Item(8, 5).Value2 = 11.7350429911083
$ws.Cells.Item(8, 7).Value2 = 3.677557044953704
$ws.Cells.Item(8, 9).Value2 = 26.90226723263567
$ws.Cells.Item(8, 11).Value2 = 10.20727077916874
$ws.Cells.Item(8, 12).Value2 = 10.21747021096918
$ws.Cells.Item(8, 13).Value2 = 14.41769488322584
$ws.Cells.Item(8, 14).Value2 = 20.90519345995563
$ws.Cells.Item(8, 15).Value2 = 27.55454481387629

$ws.Cells.Item(9, 2).Value2 = 13.11748389537095
$ws.Cells.Item(9, 3).Value2 = 7.23302687108754
$ws.Cells.Item(9, 4).Value2 = 6.258162987422088
$ws.Cells.Item(9, 5).Value2 = 11.68431274019909
$ws.Cells.Item(9, 7).Value2 = 3.67311518408727
$ws.Cells.Item(9, 9).Value2 = 26.7378050652115
$ws.Cells.Item(9, 11).Value2 = 10.57619977896125
$ws.Cells.Item(9, 12).Value2 = 10.20504053236593
$ws.Cells.Item(9, 13).Value2 = 14.49284305020299
$ws.Cells.Item(9, 14).Value2 = 20.76325000751403
$ws.Cells.Item(9, 15).Value2 = 27.43002546450089

$ws.Cells.Item(10, 2).Value2 = 13.48571397657949
$ws.Cells.Item(10, 3).Value2 = 7.371444192193659
$ws.Cells.Item(10, 4).Value2 = 6.468637313795117
$ws.Cells.Item(10, 5).Value2 = 11.65351077656414
$ws.Cells.Item(10, 7).Value2 = 3.670153979475363
$ws.Cells.Item(10, 9).Value2 = 26.63828057348887
$ws.Cells.Item(10, 11).Value2 = 10.84777279791606
$ws.Cells.Item(10, 12).Value2 = 10.20260604541099
$ws.Cells.Item(10, 13).Value2 = 14.558168459864
$ws.Cells.Item(10, 14).Value2 = 20.66734835844169
$ws.Cells.Item(10, 15).Value2 = 27.36328361611641

$ws.Cells.Item(11, 2).Value2 = 13.65228941687774
$ws.Cells.Item(11, 3).Value2 = 7.4327359075268
$ws.Cells.Item(11, 4).Value2 = 6.563761479499555
$ws.Cells.Item(11, 5).Value2 = 11.64089789271381
$ws.Cells.Item(11, 7).Value2 = 3.668871804493796
$ws.Cells.Item(11, 9).Value2 = 26.5976384668525
$ws.Cells.Item(11, 11).Value2 = 10.97089662308077
$ws.Cells.Item(11, 12).Value2 = 10.20294367399756
$ws.Cells.Item(11, 13).Value2 = 14.59001254628337
$ws.Cells.Item(11, 14).Value2 = 20.6255230801628
$ws.Cells.Item(11, 15).Value2 = 27.33830711869848

$ws.Cells.Item(12, 2).Value2 = 13.71517221907286
$ws.Cells.Item(12, 3).Value2 = 7.455696358475844
$ws.Cells.Item(12, 4).Value2 = 6.599654469616639
$ws.Cells.Item(12, 5).Value2 = 11.6363224845685
$ws.Cells.Item(12, 7).Value2 = 3.668395558592724
$ws.Cells.Item(12, 9).Value2 = 26.58291505079803
$ws.Cells.Item(12, 11).Value2 = 11.0174174148229
$ws.Cells.Item(12, 12).Value2 = 10.20327839958276
$ws.Cells.Item(12, 13).Value2 = 14.6023704394739
$ws.Cells.Item(12, 14).Value2 = 20.60994259762157
$ws.Cells.Item(12, 15).Value2 = 27.32962422802547

$ws.Cells.Item(13, 2).Value2 = 13.70163893694597
$ws.Cells.Item(13, 3).Value2 = 7.450762638575998
$ws.Cells.Item(13, 4).Value2 = 6.591930586010158
$ws.Cells.Item(13, 5).Value2 = 11.63729895518363
$ws.Cells.Item(13, 7).Value2 = 3.668497714393912
$ws.Cells.Item(13, 9).Value2 = 26.58605633091953
$ws.Cells.Item(13, 11).Value2 = 11.00740361462377
$ws.Cells.Item(13, 12).Value2 = 10.20319712537112
$ws.Cells.Item(13, 13).Value2 = 14.59969574735173
$ws.Cells.Item(13, 14).Value2 = 20.61328668338883
$ws.Cells.Item(13, 15).Value2 = 27.33145975687865

$ws.Cells.Item(14, 2).Value2 = 13.65746708766509
$ws.Cells.Item(14, 3).Value2 = 7.434629923493863
$ws.Cells.Item(14, 4).Value2 = 6.56671720862844
$ws.Cells.Item(14, 5).Value2 = 11.64051744861154
$ws.Cells.Item(14, 7).Value2 = 3.668832437611773
$ws.Cells.Item(14, 9).Value2 = 26.59641379325706
$ws.Cells.Item(14, 11).Value2 = 10.97472623659541
$ws.Cells.Item(14, 12).Value2 = 10.20296707164186
$ws.Cells.Item(14, 13).Value2 = 14.59102327993116
$ws.Cells.Item(14, 14).Value2 = 20.62423610330848
$ws.Cells.Item(14, 15).Value2 = 27.33757723245096

$ws.Cells.Item(15, 2).Value2 = 13.63038328946166
$ws.Cells.Item(15, 3).Value2 = 7.424715444089609
$ws.Cells.Item(15, 4).Value2 = 6.551255431821209
$ws.Cells.Item(15, 5).Value2 = 11.64251500994283
$ws.Cells.Item(15, 7).Value2 = 3.669038673208887
$ws.Cells.Item(15, 9).Value2 = 26.60284490608292
$ws.Cells.Item(15, 11).Value2 = 10.95469564439874
$ws.Cells.Item(15, 12).Value2 = 10.20285306909211
$ws.Cells.Item(15, 13).Value2 = 14.58574990829392
$ws.Cells.Item(15, 14).Value2 = 20.63097648486026
$ws.Cells.Item(15, 15).Value2 = 27.34142533406978

$ws.Cells.Item(16, 2).Value2 = 13.47480377971687
$ws.Cells.Item(16, 3).Value2 = 7.367404192581141
$ws.Cells.Item(16, 4).Value2 = 6.462404830467028
$ws.Cells.Item(16, 5).Value2 = 11.65436318002403
$ws.Cells.Item(16, 7).Value2 = 3.670239074523598
$ws.Cells.Item(16, 9).Value2 = 26.64102990797104
$ws.Cells.Item(16, 11).Value2 = 10.83971422592711
$ws.Cells.Item(16, 12).Value2 = 10.20261298191669
$ws.Cells.Item(16, 13).Value2 = 14.5561296170455
$ws.Cells.Item(16, 14).Value2 = 20.67011788076855
$ws.Cells.Item(16, 15).Value2 = 27.36502432084773

$ws.Cells.Item(17, 2).Value2 = 13.37907708604037
$ws.Cells.Item(17, 3).Value2 = 7.331810415045164
$ws.Cells.Item(17, 4).Value2 = 6.407710208904083
$ws.Cells.Item(17, 5).Value2 = 11.66198972840414
$ws.Cells.Item(17, 7).Value2 = 3.670992069762534
$ws.Cells.Item(17, 9).Value2 = 26.66564212915897
$ws.Cells.Item(17, 11).Value2 = 10.76903874557374
$ws.Cells.Item(17, 12).Value2 = 10.20283518299462
$ws.Cells.Item(17, 13).Value2 = 14.53849872557474
$ws.Cells.Item(17, 14).Value2 = 20.69459030359267
$ws.Cells.Item(17, 15).Value2 = 27.38088136240043

$ws.Cells.Item(18, 2).Value2 = 13.323933211102
$ws.Cells.Item(18, 3).Value2 = 7.311180644703867
$ws.Cells.Item(18, 4).Value2 = 6.376195010539579
$ws.Cells.Item(18, 5).Value2 = 11.66650802853191
$ws.Cells.Item(18, 7).Value2 = 3.671431283327648
$ws.Cells.Item(18, 9).Value2 = 26.68023437807448
$ws.Cells.Item(18, 11).Value2 = 10.72835144385584
$ws.Cells.Item(18, 12).Value2 = 10.20309906652381
$ws.Cells.Item(18, 13).Value2 = 14.52855854813875
$ws.Cells.Item(18, 14).Value2 = 20.70883576582952
$ws.Cells.Item(18, 15).Value2 = 27.39050870963795

$ws.Cells.Item(19, 2).Value2 = 13.30524985357978
$ws.Cells.Item(19, 3).Value2 = 7.304169029751921
$ws.Cells.Item(19, 4).Value2 = 6.365516059878897
$ws.Cells.Item(19, 5).Value2 = 11.66806047990997
$ws.Cells.Item(19, 7).Value2 = 3.671581044389352
$ws.Cells.Item(19, 9).Value2 = 26.68524991310171
$ws.Cells.Item(19, 11).Value2 = 10.71457048490745
$ws.Cells.Item(19, 12).Value2 = 10.20321181168365
$ws.Cells.Item(19, 13).Value2 = 14.52522763349581
$ws.Cells.Item(19, 14).Value2 = 20.71368819301333
$ws.Cells.Item(19, 15).Value2 = 27.39385537723761

$ws.Cells.Item(20, 2).Value2 = 13.38927654456003
$ws.Cells.Item(20, 3).Value2 = 7.335615762821844
$ws.Cells.Item(20, 4).Value2 = 6.413538660892113
$ws.Cells.Item(20, 5).Value2 = 11.66116424033655
$ws.Cells.Item(20, 7).Value2 = 3.670911280037917
$ws.Cells.Item(20, 9).Value2 = 26.66297699114975
$ws.Cells.Item(20, 11).Value2 = 10.77656638368905
$ws.Cells.Item(20, 12).Value2 = 10.20279745208208
$ws.Cells.Item(20, 13).Value2 = 14.54035484529683
$ws.Cells.Item(20, 14).Value2 = 20.69196763048253
$ws.Cells.Item(20, 15).Value2 = 27.37914089485175

$ws.Cells.Item(21, 2).Value2 = 13.67044721854563
$ws.Cells.Item(21, 3).Value2 = 7.439375327715288
$ws.Cells.Item(21, 4).Value2 = 6.574126776592381
$ws.Cells.Item(21, 5).Value2 = 11.63956665195374
$ws.Cells.Item(21, 7).Value2 = 3.668733869631397
$ws.Cells.Item(21, 9).Value2 = 26.59335345095505
$ws.Cells.Item(21, 11).Value2 = 10.98432752538874
$ws.Cells.Item(21, 12).Value2 = 10.20302903733481
$ws.Cells.Item(21, 13).Value2 = 14.59356252757619
$ws.Cells.Item(21, 14).Value2 = 20.62101300434918
$ws.Cells.Item(21, 15).Value2 = 27.33575933806475

$ws.Cells.Item(22, 2).Value2 = 13.85304113993483
$ws.Cells.Item(22, 3).Value2 = 7.505730835055881
$ws.Cells.Item(22, 4).Value2 = 6.678314930382746
$ws.Cells.Item(22, 5).Value2 = 11.62662171631928
$ws.Cells.Item(22, 7).Value2 = 3.667364910797875
$ws.Cells.Item(22, 9).Value2 = 26.55173778096563
$ws.Cells.Item(22, 11).Value2 = 11.11948879213105
$ws.Cells.Item(22, 12).Value2 = 10.20438571697099
$ws.Cells.Item(22, 13).Value2 = 14.63007764056293
$ws.Cells.Item(22, 14).Value2 = 20.57614244714605
$ws.Cells.Item(22, 15).Value2 = 27.31192548858829

$ws.Cells.Item(23, 2).Value2 = 13.75571376234766
$ws.Cells.Item(23, 3).Value2 = 7.470451732869939
$ws.Cells.Item(23, 4).Value2 = 6.622790167813363
$ws.Cells.Item(23, 5).Value2 = 11.6334237112137
$ws.Cells.Item(23, 7).Value2 = 3.668090614216175
$ws.Cells.Item(23, 9).Value2 = 26.57359291199747
$ws.Cells.Item(23, 11).Value2 = 11.04742170419217
$ws.Cells.Item(23, 12).Value2 = 10.20355167220371
$ws.Cells.Item(23, 13).Value2 = 14.61043183759232
$ws.Cells.Item(23, 14).Value2 = 20.59995360584947
$ws.Cells.Item(23, 15).Value2 = 27.32423239401873

$ws.Cells.Item(24, 2).Value2 = 13.38466570517593
$ws.Cells.Item(24, 3).Value2 = 7.33389588287382
$ws.Cells.Item(24, 4).Value2 = 6.410903834459177
$ws.Cells.Item(24, 5).Value2 = 11.66153702702257
$ws.Cells.Item(24, 7).Value2 = 3.670947785424319
$ws.Cells.Item(24, 9).Value2 = 26.66418052217811
$ws.Cells.Item(24, 11).Value2 = 10.77316330654079
$ws.Cells.Item(24, 12).Value2 = 10.20281408610772
$ws.Cells.Item(24, 13).Value2 = 14.53951508236982
$ws.Cells.Item(24, 14).Value2 = 20.69315279284987
$ws.Cells.Item(24, 15).Value2 = 27.37992616878642

$ws.Cells.Item(25, 2).Value2 = 12.98147065124348
$ws.Cells.Item(25, 3).Value2 = 7.180626280940136
$ws.Cells.Item(25, 4).Value2 = 6.180372460685992
$ws.Cells.Item(25, 5).Value2 = 11.69689878497644
$ws.Cells.Item(25, 7).Value2 = 3.674263523000105
$ws.Cells.Item(25, 9).Value2 = 26.77855760429257
$ws.Cells.Item(25, 11).Value2 = 10.4761182712445
$ws.Cells.Item(25, 12).Value2 = 10.20722399566942
$ws.Cells.Item(25, 13).Value2 = 14.47071447657761
$ws.Cells.Item(25, 14).Value2 = 20.80017117038467
$ws.Cells.Item(25, 15).Value2 = 27.34142533406978
